$d = $word.ActiveDocument

# --- Change the character name "Britt" to "Skylar" in the dialogue,
#     while keeping it in a run of its own as it was before. ---
$nameRange = $d.Content
[void]$nameRange.Find.Execute("Britt", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)
$nameStart = $nameRange.Start
$nameRange.Text = "Skylar"
$nameEnd = $nameStart + 6   # length of "Skylar"

# Wrapping the freshly-written text with a pair of bookmarks and then
# immediately removing them stops the run from being silently re-merged
# with its identically-formatted neighbours when the document is saved,
# so the three-run split ("Well, " / "Skylar" / ", it's like this...")
# is preserved just like in the original document.
$barrierStart = $d.Range($nameStart, $nameStart)
$barrierEnd = $d.Range($nameEnd, $nameEnd)
[void]$d.Bookmarks.Add("zzzTempBarrierStart", $barrierStart)
[void]$d.Bookmarks.Add("zzzTempBarrierEnd", $barrierEnd)
$d.Bookmarks("zzzTempBarrierStart").Delete()
$d.Bookmarks("zzzTempBarrierEnd").Delete()

# --- Move Word's automatic "_GoBack" (last edit) bookmark from its old
#     spot at the end of "He took several. " to the new last-edit
#     location. Word always drops it right where the most recent edit
#     happened; here that is inside the final paragraph, splitting
#     "Pennies are" into "Pennies ar" | "e worse ...". ---
$lastEditRange = $d.Content
[void]$lastEditRange.Find.Execute("Pennies ar", $true, $false, $false, $false, `
                                   $false, $true, 1, $false, "", 0)
$goBackRange = $d.Range($lastEditRange.End, $lastEditRange.End)
[void]$d.Bookmarks.Add("_GoBack", $goBackRange)
